$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns I and J
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Copy the existing header formatting (bold/border/centered style used by
# the other header cells, e.g. H1) onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# I0 / IF values for data rows 2..68
$i0 = @(9,7,7,9,8,8,8,9,9,7,8,8,8,9,9,8,7,8,8,7,8,9,10,9,9,10,9,7,9,5,9,8,9,9,10,9,8,9,11,9,9,9,8,9,11,8,6,9,9,9,9,9,9,7,9,9,9,8,9,9,9,6,5,5,2,6,2)
$if = @(9,7,7,9,8,8,8,9,9,7,8,8,8,9,9,8,8,8,8,7,8,9,10,9,9,10,9,7,9,5,9,8,9,9,10,9,8,9,11,9,9,9,8,9,12,8,6,9,9,9,9,9,9,8,9,9,9,9,9,9,9,6,5,5,2,6,2)

for ($idx = 0; $idx -lt $i0.Length; $idx++) {
    $r = $idx + 2
    $ws.Cells.Item($r, 9).Value = $i0[$idx]
    $ws.Cells.Item($r, 10).Value = $if[$idx]
}
